# Rename the two worksheet tabs:
#   "Enrichment using ENCORI"   -> "Network enrichment by ENCORI"
#   "Enrichment using SPONGEdb" -> "Network enrichment by SPONGEdb"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    if ($ws.Name -eq "Enrichment using ENCORI") {
        $ws.Name = "Network enrichment by ENCORI"
    }
    elseif ($ws.Name -eq "Enrichment using SPONGEdb") {
        $ws.Name = "Network enrichment by SPONGEdb"
    }
}

# Fallback (in case names already differ slightly): rename by position
# so the workbook still ends up with the two expected tab names.
if ($wb.Worksheets.Item(1).Name -ne "Network enrichment by ENCORI") {
    $wb.Worksheets.Item(1).Name = "Network enrichment by ENCORI"
}
if ($wb.Worksheets.Item(2).Name -ne "Network enrichment by SPONGEdb") {
    $wb.Worksheets.Item(2).Name = "Network enrichment by SPONGEdb"
}
